$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are formatted as text so values like "1.005" or "14.30" are not
# auto-converted to numbers by Excel, which would alter/trim their textual representation.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '25.693.34'
$ws.Range("E2").Value = '  -3.10%  '
$ws.Range("D3").Value = '1.743.86'
$ws.Range("E3").Value = '  -5.13%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '237.69'
$ws.Range("E5").Value = '  -7.93%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").Value = '0.4914'
$ws.Range("E7").Value = '  -6.01%  '
$ws.Range("D8").Value = '41.24'
$ws.Range("E8").Value = '  -8.13%  '
$ws.Range("D9").Value = '0.2431'
$ws.Range("E9").Value = '  -23.08%  '
$ws.Range("D10").Value = '0.06010'
$ws.Range("E10").Value = '  -11.26%  '
$ws.Range("D11").Value = '1.769.87'
$ws.Range("E11").Value = '  -3.69%  '
$ws.Range("D12").Value = '0.06632'
$ws.Range("E12").Value = '  -14.52%  '
$ws.Range("D13").Value = '14.30'
$ws.Range("E13").Value = '  -23.37%  '
$ws.Range("D14").Value = '0.5952'
$ws.Range("E14").Value = '  -23.38%  '
$ws.Range("D15").Value = '76.89'
$ws.Range("E15").Value = '  -12.30%  '
$ws.Range("D16").Value = '4.322'
$ws.Range("E16").Value = '  -13.52%  '
$ws.Range("D18").Value = '1.005'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '25.706.86'
$ws.Range("E19").Value = '  -3.17%  '
$ws.Range("D20").Value = '11.11'
$ws.Range("E20").Value = '  -19.69%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.000006250'
$ws.Range("E21").Value = '  -21.03%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '1.991.34'
$ws.Range("E22").Value = '  -3.87%  '
$ws.Range("D23").Value = '3.836'
$ws.Range("E23").Value = '  -16.51%  '
$ws.Range("D24").Value = '5.069'
$ws.Range("E24").Value = '  -14.88%  '
$ws.Range("D25").Value = '7.939'
$ws.Range("E25").Value = '  -14.53%  '
$ws.Range("D26").Value = '134.18'
$ws.Range("E26").Value = '  -5.91%  '
$ws.Range("D27").Value = '1.862'
$ws.Range("E27").Value = '  -15.79%  '
$ws.Range("D28").Value = '1.423'
$ws.Range("E28").Value = '  -14.86%  '
$ws.Range("D29").Value = '14.20'
$ws.Range("E29").Value = '  -15.82%  '
$ws.Range("D30").Value = '99.69'
$ws.Range("E30").Value = '  -10.70%  '
$ws.Range("D31").Value = '0.08180'
$ws.Range("E31").Value = '  -6.13%  '
$ws.Range("D32").Value = '3.593'
$ws.Range("E32").Value = '  -13.64%  '
$ws.Range("D33").Value = '1.004'
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '3.159'
$ws.Range("E34").Value = '  -22.08%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.04253'
$ws.Range("E35").Value = '  -12.70%  '
$ws.Range("D36").Value = '2.614'
$ws.Range("E36").Value = '  -8.54%  '
$ws.Range("D37").Value = '1.017'
$ws.Range("E37").Value = '  -10.18%  '
$ws.Range("D38").Value = '0.6071'
$ws.Range("E38").Value = '  -15.47%  '
$ws.Range("D39").Value = '2.690'
$ws.Range("E39").Value = '  -12.90%  '
$ws.Range("D40").Value = '2.096'
$ws.Range("E40").Value = '  -5.79%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = '101.38'
$ws.Range("E42").Value = '  -7.98%  '
$ws.Range("D43").Value = '0.01455'
$ws.Range("E43").Value = '  -15.76%  '
$ws.Range("D44").Value = '0.7856'
$ws.Range("E44").Value = '  -12.10%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3798'
$ws.Range("E45").Value = '  -20.93%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '5.134'
$ws.Range("E46").Value = '  -13.27%  '
$ws.Range("D47").Value = '6.100'
$ws.Range("E47").Value = '  -19.81%  '
$ws.Range("D48").Value = '0.05071'
$ws.Range("E48").Value = '  -12.90%  '
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '51.69'
$ws.Range("E49").Value = '  -13.21%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1032'
$ws.Range("E50").Value = '  -15.99%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '0.9994'
$ws.Range("E51").Value = '  -0.39%  '
